# Generate Report for Handoff
# Replaces the old run's GUID-based file names / timestamps with the
# new run's values across the Overview / zh-cn / de-de sheets, and
# refreshes the hyperlink display text (targets/rIds are unchanged).

$wb = $excel.ActiveWorkbook

$oldGuid = "31d5cc52-9923-4457-997f-ef273da4ed32"
$newGuid = "fa8c90cb-2def-4255-aea6-5cc38092d975"

$oldZhXlf = "31d5cc52-9923-4457-997f-ef273da4ed32.3351622299be5d19bcb1eca07fd0c671afb799de.zh-cn.xlf"
$newZhXlf = "fa8c90cb-2def-4255-aea6-5cc38092d975.638e5c72b9ee0eaa35111d7d509c9183bbb2bc4a.zh-cn.xlf"

$oldDeXlf = "31d5cc52-9923-4457-997f-ef273da4ed32.3351622299be5d19bcb1eca07fd0c671afb799de.de-de.xlf"
$newDeXlf = "fa8c90cb-2def-4255-aea6-5cc38092d975.638e5c72b9ee0eaa35111d7d509c9183bbb2bc4a.de-de.xlf"

$newHoGenerateDate = "2016-08-17 02:59:53"
$newZhHandoffDate   = "2016-08-17 02:59:45"
$newDeHandoffDate   = "2016-08-17 02:59:53"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8014483bba211f4f2daa3a032098a16274bf75e/e2e/$oldGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHoGenerateDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, $null, $null, "e2e\$newGuid.md")

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhHandoffDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, $null, $null, "$newGuid.md")

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newDeHandoffDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, $null, $null, "$newGuid.md")
